# "Menos algoritmos erroneos" - remove erroneous algorithm rows from the
# "Help.xml" sheet (sheet2): Clas-SLIQ (+ its "Decision Trees" heading),
# Clas-DMEL, Clas-GIL, Clas-DataSqueezer and Clas-Swap1.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Remove from bottom to top so row numbers used below always refer to the
# still-current (not yet renumbered) rows.
$ws2.Range("16:17").Delete()   # Clas-DataSqueezer, Clas-Swap1
$ws2.Range("8:9").Delete()     # Clas-DMEL, Clas-GIL
$ws2.Range("1:4").Delete()     # (blank), "Decision Trees" heading, Clas-SLIQ, (blank)

# Help.xml becomes the active/selected sheet, with the same selection left
# behind after the row deletions above.
$ws2.Activate()
$ws2.Range("A4:XFD5").Select()
